$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna2"
$ws.Range("C2").Value = "Epha5"
$ws.Range("D2").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 1.399915
$ws.Range("H2").Value = 4.199745
$ws.Range("I2").Value = 0.1420321708872511
$ws.Range("J2").Value = 0.1420321708872511
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1979113333333334
$ws.Range("N2").Value = 0.5937340000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2770590442033334
$ws.Range("R2").Value = 2.49353139783
$ws.Range("S2").Value = 0.1420321708872511
$ws.Range("T2").Value = 0.1420321708872511
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna2"
$ws.Range("C3").Value = "Epha5"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.864004666666666
$ws.Range("H3").Value = 17.592014
$ws.Range("I3").Value = 0.5949484882293837
$ws.Range("J3").Value = 0.5949484882293836
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1979113333333334
$ws.Range("N3").Value = 0.5937340000000001
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.160552982252889
$ws.Range("R3").Value = 10.444976840276
$ws.Range("S3").Value = 0.5949484882293837
$ws.Range("T3").Value = 0.5949484882293836
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Efna2"
$ws.Range("C4").Value = "Epha5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3047273333333333
$ws.Range("H4").Value = 0.9141819999999999
$ws.Range("I4").Value = 0.03091693758693657
$ws.Range("J4").Value = 0.03091693758693657
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1979113333333334
$ws.Range("N4").Value = 0.5937340000000001
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.06030899284311111
$ws.Range("R4").Value = 0.542780935588
$ws.Range("S4").Value = 0.03091693758693657
$ws.Range("T4").Value = 0.03091693758693657
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Efna2"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.287676333333333
$ws.Range("H5").Value = 6.863029
$ws.Range("I5").Value = 0.2321024032964287
$ws.Range("J5").Value = 0.2321024032964286
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1979113333333334
$ws.Range("N5").Value = 0.5937340000000001
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.4527570733651111
$ws.Range("R5").Value = 4.074813660286001
$ws.Range("S5").Value = 0.2321024032964287
$ws.Range("T5").Value = 0.2321024032964286